$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename table result/comment header columns (Backend -> Design)
$ws.Range('D2').Value2 = 'Ergebnis Design'
$ws.Range('E2').Value2 = 'Kommentar Design'

# Fill in Ergebnis (pass/fail 1/0) column
$ws.Range('D3').Value2 = 1
$ws.Range('D5').Value2 = 1
$ws.Range('D6').Value2 = 0
$ws.Range('D7').Value2 = 1
$ws.Range('D8').Value2 = 1
$ws.Range('D9').Value2 = 1
$ws.Range('D10').Value2 = 1
$ws.Range('D11').Value2 = 1
$ws.Range('D12').Value2 = 1
$ws.Range('D13').Value2 = 0
$ws.Range('D14').Value2 = 1
$ws.Range('D15').Value2 = 0
$ws.Range('D16').Value2 = 1
$ws.Range('D17').Value2 = 0
$ws.Range('D18').Value2 = 1
$ws.Range('D19').Value2 = 1
$ws.Range('D20').Value2 = 1
$ws.Range('D21').Value2 = 1
$ws.Range('D22').Value2 = 1
$ws.Range('D23').Value2 = 1
$ws.Range('D24').Value2 = 1
$ws.Range('D25').Value2 = 1
$ws.Range('D26').Value2 = 1
$ws.Range('D27').Value2 = 1
$ws.Range('D28').Value2 = 1
$ws.Range('D29').Value2 = 1
$ws.Range('D30').Value2 = 1
$ws.Range('D31').Value2 = 1
$ws.Range('D32').Value2 = 1

# Fill in Kommentar (free text) column
$ws.Range('E4').Value2 = 'Ich konnte den Fehlerfall nicht nachstellen.'
$ws.Range('E6').Value2 = 'Es wird mir jeden Wochentag diegleiche Route angezeigt.'
$ws.Range('E13').Value2 = 'Es wird nur die Freitagsroute angezeigt.'
$ws.Range('E15').Value2 = 'Umlaute können (teilw.) nicht richtig angezeigt werden. Bsp: beim Sausalitos "Jumbo Cocktail f??r 5???" '
$ws.Range('E17').Value2 = 'Über Button "kopieren" wird Link nicht in Zwischenspeicher kopiert. Manuelles kopieren des Llinks aus dem Textfeld ist möglich.'
$ws.Range('E33').Value2 = 'Ergebnis wird morgen nachgetragen'
